$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (A2 index 0)
$ws.Range("B2").Value = 23.1101482218129
$ws.Range("C2").Value = 4.362399951498389
$ws.Range("D2").Value = 4.072610582564339
$ws.Range("E2").Value = 10.98770551718572
$ws.Range("F2").Value = 55.77520609594337
$ws.Range("H2").Value = 7.344005520526261
$ws.Range("J2").Value = 10.54774182860854
$ws.Range("K2").Value = 19.01385598288839
$ws.Range("L2").Value = 11.57069076593213
$ws.Range("N2").Value = 26.33251820137312

# Row 3 (A3 index 1)
$ws.Range("B3").Value = 23.01054385223803
$ws.Range("C3").Value = 4.225191698887479
$ws.Range("D3").Value = 4.079697481276645
$ws.Range("E3").Value = 11.00629192790191
$ws.Range("F3").Value = 55.74029085385012
$ws.Range("H3").Value = 7.344005520526261
$ws.Range("J3").Value = 10.56528018060978
$ws.Range("K3").Value = 18.94947017374794
$ws.Range("L3").Value = 11.58520585156443
$ws.Range("N3").Value = 26.36910687975013

# Row 4 (A4 index 2)
$ws.Range("B4").Value = 22.95471465578982
$ws.Range("C4").Value = 4.140185687579641
$ws.Range("D4").Value = 4.084579972642078
$ws.Range("E4").Value = 11.01884644128625
$ws.Range("F4").Value = 55.72818307309119
$ws.Range("H4").Value = 7.344005520526261
$ws.Range("J4").Value = 10.57678869411716
$ws.Range("K4").Value = 18.91422504216142
$ws.Range("L4").Value = 11.59577601372121
$ws.Range("N4").Value = 26.39343557132002

# Row 5 (A5 index 3)
$ws.Range("B5").Value = 22.93332124749109
$ws.Range("C5").Value = 4.105417302812327
$ws.Range("D5").Value = 4.086703563336675
$ws.Range("E5").Value = 11.02425022335524
$ws.Range("F5").Value = 55.7255980272913
$ws.Range("H5").Value = 7.344005520526261
$ws.Range("J5").Value = 10.58166498311215
$ws.Range("K5").Value = 18.90095104056503
$ws.Range("L5").Value = 11.60050059881678
$ws.Range("N5").Value = 26.40381832645474

# Row 6 (A6 index 4)
$ws.Range("B6").Value = 22.92985135312875
$ws.Range("C6").Value = 4.099638247600146
$ws.Range("D6").Value = 4.08706428451957
$ws.Range("E6").Value = 11.02516490779648
$ws.Range("F6").Value = 55.72531071282231
$ws.Range("H6").Value = 7.344005520526261
$ws.Range("J6").Value = 10.58248596237753
$ws.Range("K6").Value = 18.89881293580211
$ws.Range("L6").Value = 11.60131031496328
$ws.Range("N6").Value = 26.40557068202597

# Row 7 (A7 index 5)
$ws.Range("B7").Value = 22.95442061871023
$ws.Range("C7").Value = 4.139717219775782
$ws.Range("D7").Value = 4.084608069309725
$ws.Range("E7").Value = 11.0189181530268
$ws.Range("F7").Value = 55.72813869669378
$ws.Range("H7").Value = 7.344005520526261
$ws.Range("J7").Value = 10.57685370185932
$ws.Range("K7").Value = 18.91404160353064
$ws.Range("L7").Value = 11.59583804174452
$ws.Range("N7").Value = 26.39357369926649

# Row 8 (A8 index 6)
$ws.Range("B8").Value = 23.07470985126993
$ws.Range("C8").Value = 4.315284923833465
$ws.Range("D8").Value = 4.074944146152874
$ws.Range("E8").Value = 10.99387728485627
$ws.Range("F8").Value = 55.76123281552788
$ws.Range("H8").Value = 7.344005520526261
$ws.Range("J8").Value = 10.5536357419673
$ws.Range("K8").Value = 18.99077240310485
$ws.Range("L8").Value = 11.57535166643908
$ws.Range("N8").Value = 26.344747367127

# Row 9 (A9 index 7)
$ws.Range("B9").Value = 23.35201902198728
$ws.Range("C9").Value = 4.650974767891729
$ws.Range("D9").Value = 4.06019035580608
$ws.Range("E9").Value = 10.95381647650091
$ws.Range("F9").Value = 55.89999844643978
$ws.Range("H9").Value = 7.344005520526261
$ws.Range("J9").Value = 10.51395716349031
$ws.Range("K9").Value = 19.17475127011004
$ws.Range("L9").Value = 11.54831855604587
$ws.Range("N9").Value = 26.26377667652619

# Row 10 (A10 index 8)
$ws.Range("B10").Value = 23.57975211656561
$ws.Range("C10").Value = 4.889233984693877
$ws.Range("D10").Value = 4.051886533502115
$ws.Range("E10").Value = 10.92987039635461
$ws.Range("F10").Value = 56.04670144949119
$ws.Range("H10").Value = 7.344005520526261
$ws.Range("J10").Value = 10.4883467107245
$ws.Range("K10").Value = 19.3295833339406
$ws.Range("L10").Value = 11.53644792998574
$ws.Range("N10").Value = 26.21328973538322

# Row 11 (A11 index 9)
$ws.Range("B11").Value = 23.68826327015273
$ws.Range("C11").Value = 4.995235995653541
$ws.Range("D11").Value = 4.048654557087604
$ws.Range("E11").Value = 10.92016221349583
$ws.Range("F11").Value = 56.12307334761738
$ws.Range("H11").Value = 7.344005520526261
$ws.Range("J11").Value = 10.47745935491509
$ws.Range("K11").Value = 19.4041031504389
$ws.Range("L11").Value = 11.53277721150347
$ws.Range("N11").Value = 26.19227545166539

# Row 12 (A12 index 10)
$ws.Range("B12").Value = 23.73003279284234
$ws.Range("C12").Value = 5.034992114536606
$ws.Range("D12").Value = 4.047508705461292
$ws.Range("E12").Value = 10.91665588911216
$ws.Range("F12").Value = 56.15336953863673
$ws.Range("H12").Value = 7.344005520526261
$ws.Range("J12").Value = 10.47344588919437
$ws.Range("K12").Value = 19.43289198901616
$ws.Range("L12").Value = 11.53163527125232
$ws.Range("N12").Value = 26.18459859264223

# Row 13 (A13 index 11)
$ws.Range("B13").Value = 23.72100725126482
$ws.Range("C13").Value = 5.026447645314856
$ws.Range("D13").Value = 4.047752021425639
$ws.Range("E13").Value = 10.91740348731069
$ws.Range("F13").Value = 56.14678369210614
$ws.Range("H13").Value = 7.344005520526261
$ws.Range("J13").Value = 10.47430540385768
$ws.Range("K13").Value = 19.42666674826164
$ws.Range("L13").Value = 11.53187018518095
$ws.Range("N13").Value = 26.18623945264617

# Row 14 (A14 index 12)
$ws.Range("B14").Value = 23.69168623614494
$ws.Range("C14").Value = 4.998514686352348
$ws.Range("D14").Value = 4.048558725769997
$ws.Range("E14").Value = 10.91987034255486
$ws.Range("F14").Value = 56.12553832096074
$ws.Range("H14").Value = 7.344005520526261
$ws.Range("J14").Value = 10.47712697547649
$ws.Range("K14").Value = 19.40646031093001
$ws.Range("L14").Value = 11.53267829510309
$ws.Range("N14").Value = 26.19163824413332

# Row 15 (A15 index 13)
$ws.Range("B15").Value = 23.67381380443893
$ws.Range("C15").Value = 4.981353696930713
$ws.Range("D15").Value = 4.049063004676991
$ws.Range("E15").Value = 10.92140348214262
$ws.Range("F15").Value = 56.11270377567228
$ws.Range("H15").Value = 7.344005520526261
$ws.Range("J15").Value = 10.47886949776773
$ws.Range("K15").Value = 19.3941569463227
$ws.Range("L15").Value = 11.53320557455769
$ws.Range("N15").Value = 26.19498172846403

# Row 16 (A16 index 14)
$ws.Range("B16").Value = 23.57275718296534
$ws.Range("C16").Value = 4.882254870139004
$ws.Range("D16").Value = 4.052108692134216
$ws.Range("E16").Value = 10.93052865615662
$ws.Range("F16").Value = 56.04190345508962
$ws.Range("H16").Value = 7.344005520526261
$ws.Range("J16").Value = 10.48907354582306
$ws.Range("K16").Value = 19.32479401296277
$ws.Range("L16").Value = 11.53672256912974
$ws.Range("N16").Value = 26.21470235956128

# Row 17 (A17 index 15)
$ws.Range("B17").Value = 23.51200138025469
$ws.Range("C17").Value = 4.820819368676419
$ws.Range("D17").Value = 4.054116541594304
$ws.Range("E17").Value = 10.93642984729137
$ws.Range("F17").Value = 56.00093116000879
$ws.Range("H17").Value = 7.344005520526261
$ws.Range("J17").Value = 10.49552854572491
$ws.Range("K17").Value = 19.28327628831538
$ws.Range("L17").Value = 11.53932264791674
$ws.Range("N17").Value = 26.22730042436344

# Row 18 (A18 index 16)
$ws.Range("B18").Value = 23.47752055678899
$ws.Range("C18").Value = 4.785261180453179
$ws.Range("D18").Value = 4.0553227710986
$ws.Range("E18").Value = 10.93993562243807
$ws.Range("F18").Value = 55.97827263111882
$ws.Range("H18").Value = 7.344005520526261
$ws.Range("J18").Value = 10.49931312722722
$ws.Range("K18").Value = 19.25978210564942
$ws.Range("L18").Value = 11.54098095310744
$ws.Range("N18").Value = 26.23473027253592

# Row 19 (A19 index 17)
$ws.Range("B19").Value = 23.46592651885978
$ws.Range("C19").Value = 4.773184987425664
$ws.Range("D19").Value = 4.05574001452256
$ws.Range("E19").Value = 10.94114179278251
$ws.Range("F19").Value = 55.97075701134172
$ws.Range("H19").Value = 7.344005520526261
$ws.Range("J19").Value = 10.500606870899
$ws.Range("K19").Value = 19.25189415068563
$ws.Range("L19").Value = 11.54157040602644
$ws.Range("N19").Value = 26.23727745404478

# Row 20 (A20 index 18)
$ws.Range("B20").Value = 23.51842108379953
$ws.Range("C20").Value = 4.827382587399764
$ws.Range("D20").Value = 4.053897489183736
$ws.Range("E20").Value = 10.93579011198771
$ws.Range("F20").Value = 56.0051988597822
$ws.Range("H20").Value = 7.344005520526261
$ws.Range("J20").Value = 10.49483396791938
$ws.Range("K20").Value = 19.28765612658356
$ws.Range("L20").Value = 11.53902901907159
$ws.Range("N20").Value = 26.22594031880386

# Row 21 (A21 index 19)
$ws.Range("B21").Value = 23.70028033803914
$ws.Range("C21").Value = 5.006730021582906
$ws.Range("D21").Value = 4.048319662808487
$ws.Range("E21").Value = 10.91914115849836
$ws.Range("F21").Value = 56.13174134690456
$ws.Range("H21").Value = 7.344005520526261
$ws.Range("J21").Value = 10.47629524701406
$ws.Range("K21").Value = 19.41238011335625
$ws.Range("L21").Value = 11.53243420582686
$ws.Range("N21").Value = 26.19004486726402

# Row 22 (A22 index 20)
$ws.Range("B22").Value = 23.82307808317387
$ws.Range("C22").Value = 5.121686163110411
$ws.Range("D22").Value = 4.045128840527874
$ws.Range("E22").Value = 10.90925051102115
$ws.Range("F22").Value = 56.22245801346156
$ws.Range("H22").Value = 7.344005520526261
$ws.Range("J22").Value = 10.46481626518694
$ws.Range("K22").Value = 19.49720679468923
$ws.Range("L22").Value = 11.52956979682734
$ws.Range("N22").Value = 26.16822181542914

# Row 23 (A23 index 21)
$ws.Range("B23").Value = 23.75718744640107
$ws.Range("C23").Value = 5.060551194476993
$ws.Range("D23").Value = 4.046790383085048
$ws.Range("E23").Value = 10.91443886412433
$ws.Range("F23").Value = 56.1733111255877
$ws.Range("H23").Value = 7.344005520526261
$ws.Range("J23").Value = 10.47088463781852
$ws.Range("K23").Value = 19.45163622647957
$ws.Range("L23").Value = 11.53096651889927
$ws.Range("N23").Value = 26.17971941950336

# Row 24 (A24 index 22)
$ws.Range("B24").Value = 23.51551733670428
$ws.Range("C24").Value = 4.824416094630019
$ws.Range("D24").Value = 4.0539963611273
$ws.Range("E24").Value = 10.93607898399486
$ws.Range("F24").Value = 56.00326663791232
$ws.Range("H24").Value = 7.344005520526261
$ws.Range("J24").Value = 10.49514775752357
$ws.Range("K24").Value = 19.28567483256523
$ws.Range("L24").Value = 11.53916125939789
$ws.Range("N24").Value = 26.2265546398454

# Row 25 (A25 index 23)
$ws.Range("B25").Value = 23.27268844527934
$ws.Range("C25").Value = 4.561423169918603
$ws.Range("D25").Value = 4.063734655670811
$ws.Range("E25").Value = 10.9636884997408
$ws.Range("F25").Value = 55.85457453530046
$ws.Range("H25").Value = 7.344005520526261
$ws.Range("J25").Value = 10.52406749867206
$ws.Range("K25").Value = 19.12146703894799
$ws.Range("L25").Value = 11.55422676585052
$ws.Range("N25").Value = 26.28410005372534
